# Atualização de bases das ligas, do dia: 17-06-2024 às 21:10
# Swap the data (all columns except the row-sequence column A) between each
# pair of rows so that every match's full record (id, teams, scores, odds,
# profit/loss figures, etc.) moves to the other row of the pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold per-row data (everything but "A", the sequence number).
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

# Row pairs whose contents must be swapped.
$rowPairs = @(
    @(52,53),
    @(76,77),
    @(84,85),
    @(105,106),
    @(108,109),
    @(130,131),
    @(133,134),
    @(150,151),
    @(164,165),
    @(218,219),
    @(221,222)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
